$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells (I1, J1) so they pick up the same style (bold, bordered,
# centered) used by the rest of the header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I (I0) and J (IF) column data for rows 2-12
$values = @(
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(12, 12),
    @(7, 8),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(5, 5),
    @(8, 8)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
